$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, date serial (column B), price value (columns C:G, same value in each)
$data = @(
    @(2, 43577, 26511.1),
    @(3, 43578, 26656.4),
    @(4, 43579, 26597),
    @(5, 43580, 26462.1),
    @(6, 43581, 26543.3),
    @(7, 43584, 26554.4),
    @(8, 43585, 26592.9),
    @(9, 43586, 26430.1),
    @(10, 43588, 26504.9),
    @(11, 43591, 26438.5),
    @(12, 43592, 25965.1),
    @(13, 43593, 25967.3),
    @(14, 43594, 25828.4),
    @(15, 43595, 25942.4),
    @(16, 43598, 25325),
    @(17, 43599, 25532),
    @(18, 43600, 25648),
    @(19, 43601, 25862.7),
    @(20, 43602, 25764),
    @(21, 43605, 25679.9),
    @(22, 43606, 25877.3),
    @(23, 43607, 25776.6),
    @(24, 43608, 25490.5),
    @(25, 43609, 25585.7),
    @(26, 43613, 25347.8),
    @(27, 43614, 25126.4),
    @(28, 43615, 25169.9),
    @(29, 43616, 24815),
    @(30, 43619, 24819.8),
    @(31, 43620, 25332.2),
    @(32, 43621, 25539.6),
    @(33, 43623, 25983.9),
    @(34, 43636, 26753.2),
    @(35, 43735, 26820.2),
    @(36, 43754, 27002),
    @(37, 43853, 29160.1),
    @(38, 43854, 28989.7),
    @(39, 43857, 28535.8),
    @(40, 43858, 28722.8),
    @(41, 43859, 28734.4),
    @(42, 43860, 28859.4),
    @(43, 43861, 28256),
    @(44, 43864, 28399.8),
    @(45, 43865, 28807.6),
    @(46, 43866, 29290.9),
    @(47, 43867, 29379.8),
    @(48, 43868, 29102.5),
    @(49, 43871, 29276.8),
    @(50, 43872, 29276.3),
    @(51, 43873, 29551.4),
    @(52, 43874, 29423.3),
    @(53, 43875, 29398.1),
    @(54, 43879, 29232.2),
    @(55, 43880, 29348),
    @(56, 43881, 29220),
    @(57, 43882, 28992.4),
    @(58, 43885, 27960.8),
    @(59, 43886, 27081.4),
    @(60, 43887, 26957.6),
    @(61, 43888, 25766.6),
    @(62, 43889, 25409.4),
    @(63, 43892, 26703.3),
    @(64, 43893, 25917.4),
    @(65, 43894, 27090.9),
    @(66, 43895, 26121.3),
    @(67, 43896, 25864.8),
    @(68, 43899, 23851),
    @(69, 43900, 25018.2),
    @(70, 43901, 23553.2),
    @(71, 43902, 21200.6),
    @(72, 43903, 23185.6),
    @(73, 43906, 20188.5),
    @(74, 43907, 21237.4),
    @(75, 43908, 19898.9),
    @(76, 43909, 20087.2),
    @(77, 43910, 19174),
    @(78, 43913, 18591.9),
    @(79, 43914, 20704.9),
    @(80, 43915, 21200.6),
    @(81, 43916, 22552.2),
    @(82, 43917, 21636.8),
    @(83, 43920, 22327.5),
    @(84, 43921, 21917.2),
    @(85, 43922, 20943.5),
    @(86, 43923, 21413.4),
    @(87, 43924, 21052.5),
    @(88, 43927, 22680),
    @(89, 43928, 22653.9),
    @(90, 43929, 23433.6),
    @(91, 43930, 23719.4),
    @(92, 43934, 23390.8),
    @(93, 43935, 23949.8),
    @(94, 43936, 23504.3),
    @(95, 43937, 23537.7),
    @(96, 43938, 24242.5)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $dateVal = $entry[1]
    $priceVal = $entry[2]

    $ws.Cells.Item($r, 2).Value = $dateVal   # B: Date
    $ws.Cells.Item($r, 3).Value = $priceVal  # C: Open
    $ws.Cells.Item($r, 4).Value = $priceVal  # D: High
    $ws.Cells.Item($r, 5).Value = $priceVal  # E: Low
    $ws.Cells.Item($r, 6).Value = $priceVal  # F: Close
    $ws.Cells.Item($r, 7).Value = $priceVal  # G: Adj Close
    $ws.Cells.Item($r, 8).Value = 0          # H: Volume
}

# Row 97 no longer holds data; remove it entirely since the data window shifted up
$ws.Range("B97:H97").Delete()
